# Updated CHE_grids model - 2025-08-14 22:36
#
# The "solar" worksheet has a randomized mapping (column AG, rows 4-28)
# between each technology row and the grid cell it services. This was
# regenerated upstream and the AG column values need to be re-shuffled
# to the new values below (everything else on the row - the process
# name in column S, the connections in AC/AD/AE, etc. - stays as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newGridCell = @{
    4  = "CHE_0"
    5  = "CHE_3"
    6  = "CHE_10"
    7  = "CHE_22"
    8  = "CHE_17"
    9  = "CHE_19"
    10 = "CHE_23"
    11 = "CHE_11"
    12 = "CHE_15"
    13 = "CHE_25"
    14 = "CHE_2"
    15 = "CHE_20"
    16 = "CHE_1"
    17 = "CHE_6"
    18 = "CHE_24"
    19 = "CHE_8"
    20 = "CHE_5"
    21 = "CHE_7"
    22 = "CHE_13"
    23 = "CHE_14"
    24 = "CHE_18"
    25 = "CHE_9"
    26 = "CHE_21"
    27 = "CHE_4"
    28 = "CHE_12"
}

foreach ($row in $newGridCell.Keys) {
    $ws.Cells.Item($row, 33).Value = $newGridCell[$row]
}
